# Upgrade logging system configuration
# Applies the May 2025 device-log snapshot: corrects the row-181 timestamp
# on sheet 1 and appends 8 new daily rows (182-189) to every sheet.

$wb = $excel.ActiveWorkbook

function Set-DateCell($ws, $row, $col) {
    $ws.Cells.Item($row, $col).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# --- Sheet 1: fix the row-181 timestamp (logging clock correction) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(181, 1).Value2 = 45967.46103009259

# --- Sheet 1: append rows 182-189 ---
$ws = $wb.Worksheets.Item(1)
# Row 182
$ws.Cells.Item(182, 1).Value2 = 45968.46172453704
$ws.Cells.Item(182, 2).Value2 = "0x01,0x90"
$ws.Cells.Item(182, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(182, 4).Value2 = "0x00,0xA1"
$ws.Cells.Item(182, 5).Value2 = "0x07"
$ws.Cells.Item(182, 6).Value2 = 400
$ws.Cells.Item(182, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(182, 8).Value2 = 168
$ws.Cells.Item(182, 9).Value2 = 7
Set-DateCell $ws 182 1
# Row 183
$ws.Cells.Item(183, 1).Value2 = 45969.46241898148
$ws.Cells.Item(183, 2).Value2 = "0x01,0x90"
$ws.Cells.Item(183, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(183, 4).Value2 = "0x00,0x10C"
$ws.Cells.Item(183, 5).Value2 = "0x07"
$ws.Cells.Item(183, 6).Value2 = 400
$ws.Cells.Item(183, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(183, 8).Value2 = 168
$ws.Cells.Item(183, 9).Value2 = 7
Set-DateCell $ws 183 1
# Row 184
$ws.Cells.Item(184, 1).Value2 = 45970.46311342593
$ws.Cells.Item(184, 2).Value2 = "0x01,0x90"
$ws.Cells.Item(184, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(184, 4).Value2 = "0x00,0xA2"
$ws.Cells.Item(184, 5).Value2 = "0x07"
$ws.Cells.Item(184, 6).Value2 = 400
$ws.Cells.Item(184, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(184, 8).Value2 = 164
$ws.Cells.Item(184, 9).Value2 = 7
Set-DateCell $ws 184 1
# Row 185
$ws.Cells.Item(185, 1).Value2 = 45971.46380787037
$ws.Cells.Item(185, 2).Value2 = "0x01,0x90"
$ws.Cells.Item(185, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(185, 4).Value2 = "0x00,0x11C"
$ws.Cells.Item(185, 5).Value2 = "0x07"
$ws.Cells.Item(185, 6).Value2 = 400
$ws.Cells.Item(185, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(185, 8).Value2 = 164
$ws.Cells.Item(185, 9).Value2 = 7
Set-DateCell $ws 185 1
# Row 186
$ws.Cells.Item(186, 1).Value2 = 45972.46450231481
$ws.Cells.Item(186, 2).Value2 = "0x01,0x90"
$ws.Cells.Item(186, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(186, 4).Value2 = "0x00,0xA3"
$ws.Cells.Item(186, 5).Value2 = "0x07"
$ws.Cells.Item(186, 6).Value2 = 400
$ws.Cells.Item(186, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(186, 8).Value2 = 160
$ws.Cells.Item(186, 9).Value2 = 7
Set-DateCell $ws 186 1
# Row 187
$ws.Cells.Item(187, 1).Value2 = 45973.46519675926
$ws.Cells.Item(187, 2).Value2 = "0x01,0x90"
$ws.Cells.Item(187, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(187, 4).Value2 = "0x00,0x12C"
$ws.Cells.Item(187, 5).Value2 = "0x07"
$ws.Cells.Item(187, 6).Value2 = 400
$ws.Cells.Item(187, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(187, 8).Value2 = 160
$ws.Cells.Item(187, 9).Value2 = 7
Set-DateCell $ws 187 1
# Row 188
$ws.Cells.Item(188, 1).Value2 = 45974.4658912037
$ws.Cells.Item(188, 2).Value2 = "0x01,0x90"
$ws.Cells.Item(188, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(188, 4).Value2 = "0x00,0xA4"
$ws.Cells.Item(188, 5).Value2 = "0x07"
$ws.Cells.Item(188, 6).Value2 = 400
$ws.Cells.Item(188, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(188, 8).Value2 = 156
$ws.Cells.Item(188, 9).Value2 = 7
Set-DateCell $ws 188 1
# Row 189
$ws.Cells.Item(189, 1).Value2 = 45975.46658564815
$ws.Cells.Item(189, 2).Value2 = "0x01,0x90"
$ws.Cells.Item(189, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(189, 4).Value2 = "0x00,0x13C"
$ws.Cells.Item(189, 5).Value2 = "0x07"
$ws.Cells.Item(189, 6).Value2 = 400
$ws.Cells.Item(189, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(189, 8).Value2 = 156
$ws.Cells.Item(189, 9).Value2 = 7
Set-DateCell $ws 189 1

# --- Sheet 2: append rows 182-189 ---
$ws = $wb.Worksheets.Item(2)
# Row 182
$ws.Cells.Item(182, 1).Value2 = 45968.46172453704
$ws.Cells.Item(182, 2).Value2 = "0x01,0x7c"
$ws.Cells.Item(182, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(182, 4).Value2 = "0x00,1xCC"
$ws.Cells.Item(182, 5).Value2 = "0x19"
$ws.Cells.Item(182, 6).Value2 = 380
$ws.Cells.Item(182, 7).Value2 = ("5.68432987514711e+23" -as [double])
$ws.Cells.Item(182, 8).Value2 = 200
$ws.Cells.Item(182, 9).Value2 = 25
Set-DateCell $ws 182 1
# Row 183
$ws.Cells.Item(183, 1).Value2 = 45969.46241898148
$ws.Cells.Item(183, 2).Value2 = "0x01,0x7c"
$ws.Cells.Item(183, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(183, 4).Value2 = "0x00,0xC9"
$ws.Cells.Item(183, 5).Value2 = "0x19"
$ws.Cells.Item(183, 6).Value2 = 380
$ws.Cells.Item(183, 7).Value2 = ("5.68432987514711e+23" -as [double])
$ws.Cells.Item(183, 8).Value2 = 196
$ws.Cells.Item(183, 9).Value2 = 25
Set-DateCell $ws 183 1
# Row 184
$ws.Cells.Item(184, 1).Value2 = 45970.46311342593
$ws.Cells.Item(184, 2).Value2 = "0x01,0x7c"
$ws.Cells.Item(184, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(184, 4).Value2 = "0x00,2xCC"
$ws.Cells.Item(184, 5).Value2 = "0x19"
$ws.Cells.Item(184, 6).Value2 = 380
$ws.Cells.Item(184, 7).Value2 = ("5.68432987514711e+23" -as [double])
$ws.Cells.Item(184, 8).Value2 = 192
$ws.Cells.Item(184, 9).Value2 = 25
Set-DateCell $ws 184 1
# Row 185
$ws.Cells.Item(185, 1).Value2 = 45971.46380787037
$ws.Cells.Item(185, 2).Value2 = "0x01,0x7c"
$ws.Cells.Item(185, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(185, 4).Value2 = "0x00,0xC10"
$ws.Cells.Item(185, 5).Value2 = "0x19"
$ws.Cells.Item(185, 6).Value2 = 380
$ws.Cells.Item(185, 7).Value2 = ("5.68432987514711e+23" -as [double])
$ws.Cells.Item(185, 8).Value2 = 188
$ws.Cells.Item(185, 9).Value2 = 25
Set-DateCell $ws 185 1
# Row 186
$ws.Cells.Item(186, 1).Value2 = 45972.46450231481
$ws.Cells.Item(186, 2).Value2 = "0x01,0x7c"
$ws.Cells.Item(186, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(186, 4).Value2 = "0x00,3xCC"
$ws.Cells.Item(186, 5).Value2 = "0x19"
$ws.Cells.Item(186, 6).Value2 = 380
$ws.Cells.Item(186, 7).Value2 = ("5.68432987514711e+23" -as [double])
$ws.Cells.Item(186, 8).Value2 = 188
$ws.Cells.Item(186, 9).Value2 = 25
Set-DateCell $ws 186 1
# Row 187
$ws.Cells.Item(187, 1).Value2 = 45973.46519675926
$ws.Cells.Item(187, 2).Value2 = "0x01,0x7c"
$ws.Cells.Item(187, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(187, 4).Value2 = "0x00,0xC11"
$ws.Cells.Item(187, 5).Value2 = "0x19"
$ws.Cells.Item(187, 6).Value2 = 380
$ws.Cells.Item(187, 7).Value2 = ("5.68432987514711e+23" -as [double])
$ws.Cells.Item(187, 8).Value2 = 184
$ws.Cells.Item(187, 9).Value2 = 25
Set-DateCell $ws 187 1
# Row 188
$ws.Cells.Item(188, 1).Value2 = 45974.4658912037
$ws.Cells.Item(188, 2).Value2 = "0x01,0x7c"
$ws.Cells.Item(188, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(188, 4).Value2 = "0x00,4xCC"
$ws.Cells.Item(188, 5).Value2 = "0x19"
$ws.Cells.Item(188, 6).Value2 = 380
$ws.Cells.Item(188, 7).Value2 = ("5.68432987514711e+23" -as [double])
$ws.Cells.Item(188, 8).Value2 = 184
$ws.Cells.Item(188, 9).Value2 = 25
Set-DateCell $ws 188 1
# Row 189
$ws.Cells.Item(189, 1).Value2 = 45975.46658564815
$ws.Cells.Item(189, 2).Value2 = "0x01,0x7c"
$ws.Cells.Item(189, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(189, 4).Value2 = "0x00,0xC12"
$ws.Cells.Item(189, 5).Value2 = "0x19"
$ws.Cells.Item(189, 6).Value2 = 380
$ws.Cells.Item(189, 7).Value2 = ("5.68432987514711e+23" -as [double])
$ws.Cells.Item(189, 8).Value2 = 180
$ws.Cells.Item(189, 9).Value2 = 25
Set-DateCell $ws 189 1

# --- Sheet 3: append rows 182-189 ---
$ws = $wb.Worksheets.Item(3)
# Row 182
$ws.Cells.Item(182, 1).Value2 = 45968.46172453704
$ws.Cells.Item(182, 2).Value2 = "0x00,0x6e"
$ws.Cells.Item(182, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(182, 4).Value2 = "0x00,0x5D"
$ws.Cells.Item(182, 5).Value2 = "0x15"
$ws.Cells.Item(182, 6).Value2 = 110
$ws.Cells.Item(182, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(182, 8).Value2 = 76
$ws.Cells.Item(182, 9).Value2 = 15
Set-DateCell $ws 182 1
# Row 183
$ws.Cells.Item(183, 1).Value2 = 45969.46241898148
$ws.Cells.Item(183, 2).Value2 = "0x00,0x6e"
$ws.Cells.Item(183, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(183, 4).Value2 = "0x00,0x5C"
$ws.Cells.Item(183, 5).Value2 = "0x15"
$ws.Cells.Item(183, 6).Value2 = 110
$ws.Cells.Item(183, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(183, 8).Value2 = 76
$ws.Cells.Item(183, 9).Value2 = 15
Set-DateCell $ws 183 1
# Row 184
$ws.Cells.Item(184, 1).Value2 = 45970.46311342593
$ws.Cells.Item(184, 2).Value2 = "0x00,0x6e"
$ws.Cells.Item(184, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(184, 4).Value2 = "0x00,0x6D"
$ws.Cells.Item(184, 5).Value2 = "0x15"
$ws.Cells.Item(184, 6).Value2 = 110
$ws.Cells.Item(184, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(184, 8).Value2 = 76
$ws.Cells.Item(184, 9).Value2 = 15
Set-DateCell $ws 184 1
# Row 185
$ws.Cells.Item(185, 1).Value2 = 45971.46380787037
$ws.Cells.Item(185, 2).Value2 = "0x00,0x6e"
$ws.Cells.Item(185, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(185, 4).Value2 = "0x00,0x6C"
$ws.Cells.Item(185, 5).Value2 = "0x15"
$ws.Cells.Item(185, 6).Value2 = 110
$ws.Cells.Item(185, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(185, 8).Value2 = 75
$ws.Cells.Item(185, 9).Value2 = 15
Set-DateCell $ws 185 1
# Row 186
$ws.Cells.Item(186, 1).Value2 = 45972.46450231481
$ws.Cells.Item(186, 2).Value2 = "0x00,0x6e"
$ws.Cells.Item(186, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(186, 4).Value2 = "0x00,0x7D"
$ws.Cells.Item(186, 5).Value2 = "0x15"
$ws.Cells.Item(186, 6).Value2 = 110
$ws.Cells.Item(186, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(186, 8).Value2 = 75
$ws.Cells.Item(186, 9).Value2 = 15
Set-DateCell $ws 186 1
# Row 187
$ws.Cells.Item(187, 1).Value2 = 45973.46519675926
$ws.Cells.Item(187, 2).Value2 = "0x00,0x6e"
$ws.Cells.Item(187, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(187, 4).Value2 = "0x00,0x7C"
$ws.Cells.Item(187, 5).Value2 = "0x15"
$ws.Cells.Item(187, 6).Value2 = 110
$ws.Cells.Item(187, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(187, 8).Value2 = 75
$ws.Cells.Item(187, 9).Value2 = 15
Set-DateCell $ws 187 1
# Row 188
$ws.Cells.Item(188, 1).Value2 = 45974.4658912037
$ws.Cells.Item(188, 2).Value2 = "0x00,0x6e"
$ws.Cells.Item(188, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(188, 4).Value2 = "0x00,0x8D"
$ws.Cells.Item(188, 5).Value2 = "0x15"
$ws.Cells.Item(188, 6).Value2 = 110
$ws.Cells.Item(188, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(188, 8).Value2 = 74
$ws.Cells.Item(188, 9).Value2 = 15
Set-DateCell $ws 188 1
# Row 189
$ws.Cells.Item(189, 1).Value2 = 45975.46658564815
$ws.Cells.Item(189, 2).Value2 = "0x00,0x6e"
$ws.Cells.Item(189, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(189, 4).Value2 = "0x00,0x8C"
$ws.Cells.Item(189, 5).Value2 = "0x15"
$ws.Cells.Item(189, 6).Value2 = 110
$ws.Cells.Item(189, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(189, 8).Value2 = 74
$ws.Cells.Item(189, 9).Value2 = 15
Set-DateCell $ws 189 1

# --- Sheet 4: append rows 182-189 ---
$ws = $wb.Worksheets.Item(4)
# Row 182
$ws.Cells.Item(182, 1).Value2 = 45968.46172453704
$ws.Cells.Item(182, 2).Value2 = "0x00,0x82"
$ws.Cells.Item(182, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(182, 4).Value2 = "0x00,0x64"
$ws.Cells.Item(182, 5).Value2 = "0x9"
$ws.Cells.Item(182, 6).Value2 = 130
$ws.Cells.Item(182, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(182, 8).Value2 = 100
$ws.Cells.Item(182, 9).Value2 = 9
Set-DateCell $ws 182 1
# Row 183
$ws.Cells.Item(183, 1).Value2 = 45969.46241898148
$ws.Cells.Item(183, 2).Value2 = "0x00,0x82"
$ws.Cells.Item(183, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(183, 4).Value2 = "0x00,0x63"
$ws.Cells.Item(183, 5).Value2 = "0x9"
$ws.Cells.Item(183, 6).Value2 = 130
$ws.Cells.Item(183, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(183, 8).Value2 = 96
$ws.Cells.Item(183, 9).Value2 = 9
Set-DateCell $ws 183 1
# Row 184
$ws.Cells.Item(184, 1).Value2 = 45970.46311342593
$ws.Cells.Item(184, 2).Value2 = "0x00,0x82"
$ws.Cells.Item(184, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(184, 4).Value2 = "0x00,0x62"
$ws.Cells.Item(184, 5).Value2 = "0x9"
$ws.Cells.Item(184, 6).Value2 = 130
$ws.Cells.Item(184, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(184, 8).Value2 = 95
$ws.Cells.Item(184, 9).Value2 = 9
Set-DateCell $ws 184 1
# Row 185
$ws.Cells.Item(185, 1).Value2 = 45971.46380787037
$ws.Cells.Item(185, 2).Value2 = "0x00,0x82"
$ws.Cells.Item(185, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(185, 4).Value2 = "0x00,0x61"
$ws.Cells.Item(185, 5).Value2 = "0x9"
$ws.Cells.Item(185, 6).Value2 = 130
$ws.Cells.Item(185, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(185, 8).Value2 = 93
$ws.Cells.Item(185, 9).Value2 = 9
Set-DateCell $ws 185 1
# Row 186
$ws.Cells.Item(186, 1).Value2 = 45972.46450231481
$ws.Cells.Item(186, 2).Value2 = "0x00,0x82"
$ws.Cells.Item(186, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(186, 4).Value2 = "0x00,0x60"
$ws.Cells.Item(186, 5).Value2 = "0x9"
$ws.Cells.Item(186, 6).Value2 = 130
$ws.Cells.Item(186, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(186, 8).Value2 = 92
$ws.Cells.Item(186, 9).Value2 = 9
Set-DateCell $ws 186 1
# Row 187
$ws.Cells.Item(187, 1).Value2 = 45973.46519675926
$ws.Cells.Item(187, 2).Value2 = "0x00,0x82"
$ws.Cells.Item(187, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(187, 4).Value2 = "0x00,0x59"
$ws.Cells.Item(187, 5).Value2 = "0x9"
$ws.Cells.Item(187, 6).Value2 = 130
$ws.Cells.Item(187, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(187, 8).Value2 = 91
$ws.Cells.Item(187, 9).Value2 = 9
Set-DateCell $ws 187 1
# Row 188
$ws.Cells.Item(188, 1).Value2 = 45974.4658912037
$ws.Cells.Item(188, 2).Value2 = "0x00,0x82"
$ws.Cells.Item(188, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(188, 4).Value2 = "0x00,0x58"
$ws.Cells.Item(188, 5).Value2 = "0x9"
$ws.Cells.Item(188, 6).Value2 = 130
$ws.Cells.Item(188, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(188, 8).Value2 = 91
$ws.Cells.Item(188, 9).Value2 = 9
Set-DateCell $ws 188 1
# Row 189
$ws.Cells.Item(189, 1).Value2 = 45975.46658564815
$ws.Cells.Item(189, 2).Value2 = "0x00,0x82"
$ws.Cells.Item(189, 3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(189, 4).Value2 = "0x00,0x57"
$ws.Cells.Item(189, 5).Value2 = "0x9"
$ws.Cells.Item(189, 6).Value2 = 130
$ws.Cells.Item(189, 7).Value2 = ("5.68631262647113e+23" -as [double])
$ws.Cells.Item(189, 8).Value2 = 90
$ws.Cells.Item(189, 9).Value2 = 9
Set-DateCell $ws 189 1
